# Populate Sheet1 column A with the values A, B, C (stored as shared
# strings) and leave the selection/active cell on A3, matching the
# author's edit (git-ignore housekeeping "and some change also").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "A"
$ws.Range("A2").Value = "B"
$ws.Range("A3").Value = "C"

# Match the saved selection/active cell (A3) from the diff.
$ws.Range("A3").Select() | Out-Null
